$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (old C..F shift to D..G)
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(3).ColumnWidth = 13.333333333333334

# New header cell
$ws.Cells.Item(1,3).Value = "QueryParameters"

# New row 6 data (insert shared strings in the required order)
$ws.Cells.Item(6,3).Value = "name=happy, priya=nothing"
$ws.Cells.Item(6,2).Value = "name=new, priya=nothing"
$ws.Cells.Item(6,1).Value = "users/1/{name}/{priya}"
$ws.Cells.Item(6,5).Value = "Get"

# Copy formatting from existing rows/cells so the new row matches style
$ws.Cells.Item(5,1).Copy()
$ws.Cells.Item(6,1).PasteSpecial(-4122)

$ws.Cells.Item(1,2).Copy()
$ws.Cells.Item(6,2).PasteSpecial(-4122)
$ws.Cells.Item(6,3).PasteSpecial(-4122)

$ws.Cells.Item(5,5).Copy()
$ws.Cells.Item(6,5).PasteSpecial(-4122)

Write-Host ("Dimension: " + $ws.UsedRange.Address())
